$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 5 (previously an empty gap between the last data row and the
# pre-formatted blank row 6) with a new review entry. No rows are shifted;
# row 6 below stays exactly where it is.
$ws.Range("A5").Value = "com.hamxa.shaynachim"
$ws.Range("B5").Value = "bitcoin"
$ws.Range("C5").Value = "milleradir327@gmail.com "
$ws.Range("D5").Value = "itamaramir2@gmail.com"
$ws.Range("E5").Value = "27/5/2019 15:59"
$ws.Range("F5").Value = "I found this game really spontaneous and it is great"
$ws.Range("G5").Value = "no"

# Match formatting used by the other data rows (rows 2-4), including row height.
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# Update the selection to match the new active cell/selection used after the edit.
$ws.Range("C5:D5").Select()
$ws.Range("D5").Activate()
